# Regenerate s_vals data to filter save games.
# Updates the numeric stat columns (B:E) and the computed sum column (G)
# for each data row (rows 2-6) on the active sheet. Column A (date labels)
# and column F (Win flag) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ B = 3.286832544864788;    C = 1.655778082260271;   D = 0.7527432677738641;  E = 10.19245300693656; G = 15.88780690183548 }
    3 = @{ B = 0.04271373187048222;  C = 0.306821227259698;   D = 3.537761648806719;   E = 10.19245300693656; G = 14.07974961487346 }
    4 = @{ B = 0.1190320826869504;   C = 0.306821227259698;   D = 3.537761648806719;   E = 0.4942365360607697; G = 4.457851494814137 }
    5 = @{ B = 0.0006408296065709695;C = 0.04071648406533734; D = 3.537761648806719;   E = 10.19245300693656; G = 13.77157196941518 }
    6 = @{ B = 3.286832544864788;    C = 1.655778082260271;   D = 3.537761648806719;   E = 10.19245300693656; G = 18.67282528286833 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
